$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase de Grupos")

# Jogo (match) results for "Grupo C":
#  Row 19: Benfica 1 x 0 Bayern de Munique
#  Row 20: Auckland City 1 x 1 Boca Juniors
# NOTE: Deliberately write the "H" (away-score) cells before the "F"
# (home-score) cells on each row. Several dependent formulas in this
# sheet are nested IF(...="",...) chains; writing the outer-tested cell
# (F) last ensures the dependency ripple lands on every dependent cell
# (K/L/T/etc.) once both scores are present.
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 1

# Restore the active sheet/selection to where the user left off (F25).
$ws.Activate() | Out-Null
$ws.Range("F25").Select() | Out-Null
